# The workbook was updated from the 3rd-trimester 2021 report to the
# 4th-trimester 2021 report: the reporting period dates (B8:C8) and the
# "last updated" dates (AJ8:AK8) on sheet "Reporte de Formatos" move
# forward one quarter, and the last user selection on that sheet moves
# to AL8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Periodo que se informa: 01/10/2021 - 31/12/2021 (previously 01/07/2021 - 30/09/2021)
$ws.Range("B8").Value = 44470
$ws.Range("C8").Value = 44561

# Fecha de validación / actualización: 10/01/2022 (previously 11/10/2021)
$ws.Range("AJ8").Value = 44571
$ws.Range("AK8").Value = 44571

# Reflect the author's last active cell/selection on this sheet.
$ws.Activate() | Out-Null
$ws.Range("AL8").Select() | Out-Null
